{"js": "// 1. Remove the stray \"_GoBack\" bookmark that currently sits after \"longest\".\n//    Word re-creates this bookmark at the last edit position on every save,\n//    so it needs to move from its old spot to the new one created below.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2. Expand \"...his own status. In this case...\" into\n//    \"...his own status by self-reporting. In this case...\".\nconst body = context.document.body;\nconst target = body.search(\"his own status. In this case\", { matchCase: true });\ntarget.load(\"items\");\nawait context.sync();\n\nif (target.items.length > 0) {\n  target.items[0].insertText(\n    \"his own status by self-reporting. In this case\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// 3. Re-create the \"_GoBack\" bookmark as a zero-length mark sitting right\n//    between \"self-report\" and \"ing\" (i.e. immediately after \"self-report\"),\n//    matching where Word last left the editing cursor.\nconst selfReport = body.search(\"self-report\", { matchCase: true });\nselfReport.load(\"items\");\nawait context.sync();\n\nif (selfReport.items.length > 0) {\n  const insertionPoint = selfReport.items[0].getRange(Word.RangeLocation.end);\n  insertionPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Remove the stray \"_GoBack\" bookmark that currently sits after \"longest\".\n#    Word re-creates this bookmark at the last edit position on every save,\n#    so it needs to move from its old spot to the new one created below.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2. Expand \"...his own status. In this case...\" into\n#    \"...his own status by self-reporting. In this case...\".\n$range = $d.Content\n$find = $range.Find\n$find.Text = \"own status. In this case\"\n$found = $find.Execute()\nif ($found) {\n    $range.Text = \"own status by self-reporting. In this case\"\n}\n\n# 3. Re-create the \"_GoBack\" bookmark as a zero-length mark sitting right\n#    between \"self-report\" and \"ing\" (i.e. immediately after \"self-report\"),\n#    matching where Word last left the editing cursor.\n$range2 = $d.Content\n$find2 = $range2.Find\n$find2.Text = \"self-report\"\n$found2 = $find2.Execute()\nif ($found2) {\n    $range2.Collapse(0)  # wdCollapseEnd\n    $d.Bookmarks.Add(\"_GoBack\", $range2)\n}\n"}
